# "refactor, meet initial requirements"
#
# Adds a "percentage" column (hours / total-hours * 100) to both
# worksheets:
#   - "PI hours":   name | hours | percentage | dept      (was: name | hours | dept)
#   - "dept hours": dept | hours | percentage             (was: dept | hours)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "PI hours"
# A "percentage" column is inserted between "hours" and "dept", so the
# existing "dept" column (D) shifts right to E. Inserting a whole column
# (rather than just writing into D) keeps the existing data/format for
# the old column D attached to the cells that move to E, and Excel
# automatically extends the header's bordered/bold/centered style (and
# the index column's style) into the new column D.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("PI hours")

$ws1.Columns.Item(4).Insert()

$ws1.Range("D1").Value = "percentage"
$ws1.Range("D2").Value = 72.72727272727273   # 8 / (8 + 3) * 100
$ws1.Range("D3").Value = 27.27272727272727   # 3 / (8 + 3) * 100

# ---------------------------------------------------------------------
# Sheet 2: "dept hours"
# "percentage" is appended as a new trailing column D. Copy the existing
# header style from C1 (hours) onto D1 first via paste-special-formats so
# the new header cell reuses the workbook's existing bordered/bold/
# centered style instead of minting a duplicate one, then fill in the
# header text and the computed values.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("dept hours")

$ws2.Range("C1").Copy() | Out-Null
$ws2.Range("D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws2.Range("D1").Value = "percentage"
$ws2.Range("D2").Value = 40.74074074074074   # 11 / (11 + 8 + 8) * 100
$ws2.Range("D3").Value = 29.62962962962963   # 8  / (11 + 8 + 8) * 100
$ws2.Range("D4").Value = 29.62962962962963   # 8  / (11 + 8 + 8) * 100

Write-Host "edit complete"
